{"js": "// Replace the \"Post-Conditions\" cell content \"Kh\u00f4ng c\u00f3\" with the full\n// sentence describing who can view invoice details.\nconst body = context.document.body;\n\nconst searchResults = body.search(\"Kh\u00f4ng c\u00f3\", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items/text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find text \"Kh\u00f4ng c\u00f3\" in the document.');\n}\n\n// There is exactly one occurrence in this document (the Post-Conditions\n// table cell), but guard against future duplicates by only touching the\n// first match.\nsearchResults.items[0].insertText(\n  \"Qu\u1ea3n tr\u1ecb ho\u1eb7c nh\u00e2n vi\u00ean xem \u0111\u01b0\u1ee3c th\u00f4ng tin chi ti\u1ebft h\u00f3a \u0111\u01a1n.\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# Replace the \"Post-Conditions\" cell content \"Kh\u00f4ng c\u00f3\" with the full\n# sentence describing who can view invoice details.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Kh\u00f4ng c\u00f3\"\n$find.Replacement.Text = \"Qu\u1ea3n tr\u1ecb ho\u1eb7c nh\u00e2n vi\u00ean xem \u0111\u01b0\u1ee3c th\u00f4ng tin chi ti\u1ebft h\u00f3a \u0111\u01a1n.\"\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)  # wdReplaceAll\n"}
